$d = $word.ActiveDocument

# The three image placeholders (two "Pearson logo" ones in the first-page /
# default footers, one "BTec logo" one in the first-page header) had their
# display Name swapped between image1/image2. InlineShape does not expose a
# settable Name in the Word OM, and the underlying picture name lives in both
# wp:docPr/@name and pic:cNvPr/@name, so the safest route is a targeted,
# literal (non-regex) text substitution over the package-level OOXML.
$xml = $d.WordOpenXML

# BTec logo (header, first page) - image1.jpg -> image2.jpg
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>')
$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>')

# Pearson logo (footer, default) id="2" - image2.png -> image1.png
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>')

# Pearson logo (footer, first page) id="3" - image2.png -> image1.png
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/>')

# Both Pearson logo pic:cNvPr entries (inside nvPicPr) share identical text
# and both need the same rename, so one global replace covers both.
$xml = $xml.Replace(
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>',
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>')

$d.WordOpenXML = $xml
